$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '45.445.26'
$ws.Range("E2").Value = '  +6.39%  '
$ws.Range("D3").Value = '2.382.69'
$ws.Range("E3").Value = '  +4.57%  '
$ws.Range("E4").Value = '  -0.77%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.64'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.39%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '111.36'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +7.18%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.637'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.71%  '
$ws.Range("E8").Value = '  -0.31%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.630'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.09'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +8.86%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0931'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.91%  '
$ws.Range("E12").Value = '  +5.92%  '
$ws.Range("E13").Value = '  +4.46%  '
$ws.Range("E14").Value = '  +0.89%  '
$ws.Range("E15").Value = '  +4.88%  '
$ws.Range("D16").Value = '2.745.79'
$ws.Range("E16").Value = '  +4.72%  '
$ws.Range("D17").Value = '2.383.84'
$ws.Range("E17").Value = '  +4.75%  '
$ws.Range("D18").Value = '45.442.87'
$ws.Range("E18").Value = '  +7.27%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.64'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.82%  '
$ws.Range("E20").Value = '  +3.96%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.19'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.23%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '75.15'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.29%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.53'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.79%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '269.56'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.76%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.31'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +6.72%  '
$ws.Range("E26").Value = '  -0.69%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.31'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +6.58%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.56'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +8.65%  '
$ws.Range("E29").Value = '  +0.50%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.93'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.60%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '38.62'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +8.47%  '
$ws.Range("E32").Value = '  +11.68%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '169.80'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.53%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.02'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +18.53%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.133'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.14%  '
$ws.Range("E36").Value = '  +6.46%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.85'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +8.31%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.08'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +13.62%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0365'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.08%  '
$ws.Range("E40").Value = '  +6.70%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.73'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +12.49%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '105.53'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +7.20%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '13.71'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +15.42%  '
$ws.Range("E44").Value = '  +6.99%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '71.25'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.53%  '
$ws.Range("E46").Value = '  +0.37%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '117.87'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +7.60%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.82'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +13.46%  '
$ws.Range("E49").Value = '  +21.01%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.29'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +8.17%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '79.26'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.58%  '
